$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 633, shifting existing rows 633:687 down to 634:688
$ws.Rows.Item(633).Insert()

# Populate the newly inserted row 633 with the new data record
$ws.Cells.Item(633, 1).Value = 5
$ws.Cells.Item(633, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(633, 3).Value = "Maule"
$ws.Cells.Item(633, 4).Value = 45106
$ws.Cells.Item(633, 5).Value = 7
$ws.Cells.Item(633, 6).Value = 100112043
$ws.Cells.Item(633, 7).Value = "Pepino ensalada"
$ws.Cells.Item(633, 8).Value = "Sin especificar"
$ws.Cells.Item(633, 9).Value = "Primera"
$ws.Cells.Item(633, 10).Value = 300
$ws.Cells.Item(633, 11).Value = 15000
$ws.Cells.Item(633, 12).Value = 15000
$ws.Cells.Item(633, 13).Value = 15000
$ws.Cells.Item(633, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(633, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(633, 16).Value = 250
$ws.Cells.Item(633, 17).Value = 60
$ws.Cells.Item(633, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D in the surrounding rows
$ws.Cells.Item(633, 4).NumberFormat = $ws.Cells.Item(634, 4).NumberFormat
